{"js": "// Replace the first \"Ingr\u00e9dients\" bullet text and add three new bullet\n// items right after it (before the existing blank spacer paragraph),\n// matching the recipe update to \"Mac & Cheese Courge Musqu\u00e9e\":\n//   \"2kg courge musqu\u00e9e, pel\u00e9e\" -> \"2000g macaronis\"\n//   + \"1500g fromage cheddar orange\"\n//   + \"2.27kg de courge musqu\u00e9e\"\n//   + \"Cuire les morceaux de courges minimum 1h dans l'eau bouillante\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the ingredient bullet that currently reads\n// \"2kg courge musqu\u00e9e, pel\u00e9e\" (first list item under \"Ingr\u00e9dients\").\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"2kg courge musqu\u00e9e, pel\u00e9e\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find paragraph \"2kg courge musqu\u00e9e, pel\u00e9e\"');\n}\n\n// Swap its text for the new macaroni quantity (keeps the paragraph's\n// list formatting / run properties intact).\ntarget.insertText(\"2000g macaronis\", Word.InsertLocation.replace);\n\n// Insert the three new ingredient lines right after it, each copying the\n// anchor paragraph's list style (bullet list \"Paragraphedeliste\", numId 3).\nconst cheeseP = target.insertParagraph(\"1500g fromage cheddar orange\", Word.InsertLocation.after);\nconst squashP = cheeseP.insertParagraph(\"2.27kg de courge musqu\u00e9e\", Word.InsertLocation.after);\n\n// The last new line (cooking instructions) keeps the bullet list style but,\n// unlike the two ingredient lines above it, has no hanging indent override\n// in the source edit. Build it as a fresh paragraph + re-attach to the same\n// list instead of cloning, so no <w:ind> survives.\nconst cookP = squashP.insertParagraph(\n  \"Cuire les morceaux de courges minimum 1h dans l\\u2019eau bouillante\",\n  Word.InsertLocation.after\n);\ncookP.style = \"Paragraphedeliste\";\ncookP.attachToList(3, 0);\ncookP.spaceAfter = 0;\n\nawait context.sync();\n", "ps1": "# Replace the first \"Ingr\u00e9dients\" bullet text and add three new bullet\n# items right after it (before the existing blank spacer paragraph),\n# matching the recipe update to \"Mac & Cheese Courge Musqu\u00e9e\":\n#   \"2kg courge musqu\u00e9e, pel\u00e9e\" -> \"2000g macaronis\"\n#   + \"1500g fromage cheddar orange\"\n#   + \"2.27kg de courge musqu\u00e9e\"\n#   + \"Cuire les morceaux de courges minimum 1h dans l'eau bouillante\"\n\n$d = $word.ActiveDocument\n\n# Locate the ingredient bullet that currently reads\n# \"2kg courge musqu\u00e9e, pel\u00e9e\" (first list item under \"Ingr\u00e9dients\").\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    if ($candidate.Range.Text.TrimEnd(\"`r\", \"`n\") -eq \"2kg courge musqu\u00e9e, pel\u00e9e\") {\n        $target = $candidate\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not find paragraph '2kg courge musqu\u00e9e, pel\u00e9e'\"\n}\n\n# Remember the bullet list's template so the new paragraphs can rejoin the\n# same list (same numId) instead of minting a brand-new one.\n$listTemplate = $target.Range.ListFormat.ListTemplate\n\n$targetIndex = $target.Index\n\n# Swap its text for the new macaroni quantity (keeps the paragraph's list\n# formatting / run properties intact). No trailing carriage return here -\n# that would insert an extra empty paragraph rather than just replacing\n# this paragraph's text in place.\n$target.Range.Text = \"2000g macaronis\"\n\n# Re-fetch by index to stay safe across the text-range swap above.\n$target = $d.Paragraphs.Item($targetIndex)\n\n# Insert the two ingredient lines right after it; InsertParagraphAfter\n# clones the anchor paragraph's list formatting (style + numId=3 + the\n# 709/-283 hanging indent), which matches these two lines exactly.\n$target.Range.InsertParagraphAfter()\n$cheeseP = $d.Paragraphs.Item($target.Index + 1)\n$cheeseP.Range.Text = \"1500g fromage cheddar orange\"\n\n$cheeseP.Range.InsertParagraphAfter()\n$squashP = $d.Paragraphs.Item($cheeseP.Index + 1)\n$squashP.Range.Text = \"2.27kg de courge musqu\u00e9e\"\n\n# The last new line (cooking instructions) keeps the bullet list style but,\n# unlike the two ingredient lines above it, has no hanging indent override\n# in the source edit. Reset its paragraph style (which also drops the\n# cloned indent) then rejoin the same list (numId=3) explicitly.\n$squashP.Range.InsertParagraphAfter()\n$cookP = $d.Paragraphs.Item($squashP.Index + 1)\n$cookP.Range.Text = \"Cuire les morceaux de courges minimum 1h dans l\" + [char]0x2019 + \"eau bouillante\"\n$cookP.Style = \"Paragraphedeliste\"\n$cookP.SpaceAfter = 0\n$cookP.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)\n"}
